# Clean up madness spells and add them to spell list
$wb = $excel.ActiveWorkbook

# --- Add the three new "madness" spells to the Spells sheet ---
$spells = $wb.Worksheets.Item("Spells")

# Row 48: Creeping Smite
$spells.Cells.Item(48, 1).Value = "Creeping Smite"
$spells.Cells.Item(48, 2).Value = "New"
$spells.Cells.Item(48, 3).Value = 1
$spells.Cells.Item(48, 4).Value = "Evocation"
$spells.Cells.Item(48, 5).Value = "No"
$spells.Cells.Item(48, 6).Value = "No"
$spells.Cells.Item(48, 7).Value = "No"
$spells.Cells.Item(48, 8).Value = "No"
$spells.Cells.Item(48, 9).Value = "Yes"
$spells.Cells.Item(48, 10).Value = "No"
$spells.Cells.Item(48, 11).Value = "No"
$spells.Cells.Item(48, 12).Value = "No"
$spells.Cells.Item(48, 13).Value = "No"
$spells.Cells.Item(48, 14).Value = "Playtest Ready"
$spells.Cells.Item(48, 15).Value = "Not Released"

# Row 49: Fiery Temper
$spells.Cells.Item(49, 1).Value = "Fiery Temper"
$spells.Cells.Item(49, 2).Value = "New"
$spells.Cells.Item(49, 3).Value = 1
$spells.Cells.Item(49, 4).Value = "Evocation"
$spells.Cells.Item(49, 5).Value = "No"
$spells.Cells.Item(49, 6).Value = "Yes"
$spells.Cells.Item(49, 7).Value = "No"
$spells.Cells.Item(49, 8).Value = "No"
$spells.Cells.Item(49, 9).Value = "No"
$spells.Cells.Item(49, 10).Value = "No"
$spells.Cells.Item(49, 11).Value = "Yes"
$spells.Cells.Item(49, 12).Value = "Yes"
$spells.Cells.Item(49, 13).Value = "Yes"
$spells.Cells.Item(49, 14).Value = "Playtest Ready"
$spells.Cells.Item(49, 15).Value = "Not Released"

# Row 50: Senseless Rage
$spells.Cells.Item(50, 1).Value = "Senseless Rage"
$spells.Cells.Item(50, 2).Value = "New"
$spells.Cells.Item(50, 3).Value = 1
$spells.Cells.Item(50, 4).Value = "Enchantment"
$spells.Cells.Item(50, 5).Value = "No"
$spells.Cells.Item(50, 6).Value = "No"
$spells.Cells.Item(50, 7).Value = "Yes"
$spells.Cells.Item(50, 8).Value = "No"
$spells.Cells.Item(50, 9).Value = "No"
$spells.Cells.Item(50, 10).Value = "No"
$spells.Cells.Item(50, 11).Value = "Yes"
$spells.Cells.Item(50, 12).Value = "Yes"
$spells.Cells.Item(50, 13).Value = "No"
$spells.Cells.Item(50, 14).Value = "Playtest Ready"
$spells.Cells.Item(50, 15).Value = "Not Released"

# Update the Spells sheet view: move selection below the newly-added rows
$spells.Activate()
[void]$spells.Range("A51").Select()

# --- Add the "Madness Spells" document entry to DMResources ---
$dm = $wb.Worksheets.Item("DMResources")
$dm.Cells.Item(5, 1).Value = "Madness Spells"
$dm.Cells.Item(5, 2).Value = "Madness Spells"
$dm.Cells.Item(5, 3).Value = "Playtest Ready"
$dm.Cells.Item(5, 4).Value = "Not Released"
$dm.Hyperlinks.Add($dm.Range("B5"), "https://editor.gmbinder.com/documents/edit/-PLACEHOLDER12345678")

# Make DMResources the active/selected sheet (tab moved here)
$dm.Activate()
[void]$dm.Range("D8").Select()

# --- Subclasses sheet view tweak (tab selection moved away, scroll position changed) ---
$subclasses = $wb.Worksheets.Item("Subclasses")
$subclasses.Activate()
[void]$subclasses.Range("E80").Select()

# Re-activate DMResources last so it ends up the selected/active tab
$dm.Activate()
